$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 19 for new Movement Phase events (e018-e021)
$ws.Rows.Item(19).Resize(4).Insert()

# Set the event-code column (A) for the new rows
$ws.Range("A19").Value = "e018"
$ws.Range("A20").Value = "e019"
$ws.Range("A21").Value = "e020"
$ws.Range("A22").Value = "e021"

# Long descriptive text for each new event (column B)
$e020text = @'
<Bold>e020 Enemy Strenth Check</Bold> 
<InlineUIContainer><Button Content='r4.53' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Check any one adjacent area to your task force for estimating enemy strength. Click on one of the adjacent regions highlighted blue.
<LineBreak/><LineBreak/>
Roll 1D and consult the <InlineUIContainer><Button Content='r4.51' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table. The area is marked with a Light, Medium, or Heavy marker.
<LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer> 
<LineBreak/><LineBreak/>
'@
$e019text = @'
<Bold>e019 Set Exit Area</Bold> 
<InlineUIContainer><Button Content='r4.52' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
On the <InlineUIContainer><Button Content='Exit Areas' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table, roll 1D and cross reference the number with the Start Area marker 
<InlineUIContainer><Button Content='r4.51' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.  
<LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer> 
<LineBreak/><LineBreak/>
'@
$e018text = @'
<Bold>e018 Set Start Area</Bold> 
<InlineUIContainer><Button Content='r4.51' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Around the edge of the Movement Board, there are 10 areas number 1-10. The area is marked with the Start Area and Task Force markers. 
<LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$e021text = @'
<Bold>e021 Choose Operations</Bold> 
<InlineUIContainer><Button Content='r4.54' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Choose one of following options:<LineBreak/>
 <InlineUIContainer><Button Content='Additional' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Area Check<LineBreak/>
 <InlineUIContainer><Button Content='Artillery' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Support <LineBreak/>
 <InlineUIContainer><Button Content='Air' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Strike  <LineBreak/>
 <InlineUIContainer><Button Content='Attempt' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Resupply <LineBreak/>
Alternatively, click on an adjacent highlighted area on the Movement Board to enter that area.
<LineBreak/><LineBreak/>
'@

$ws.Range("B21").Value = $e020text
$ws.Range("B20").Value = $e019text
$ws.Range("B19").Value = $e018text
$ws.Range("B22").Value = $e021text

# Row heights for the new Movement Phase event rows
$ws.Rows.Item(19).RowHeight = 105
$ws.Rows.Item(20).RowHeight = 135
$ws.Rows.Item(21).RowHeight = 150
$ws.Rows.Item(22).RowHeight = 150

# Minor row-height tweaks to existing rows (text reflowed slightly wider)
$ws.Rows.Item(12).RowHeight = 105
$ws.Rows.Item(13).RowHeight = 120

# Leave the selection on B20, matching where editing left off
$ws.Range("B20").Select()
